$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TELEXo")

# Update quantities in row 24 and row 25 (1k trigger out resistor reversion)
$ws.Range("A24").Value = 8
$ws.Range("A25").Value = 8

# Move the "R102,R103,R104,R105" label from G25 to G24
$ws.Range("G24").Value = "R102,R103,R104,R105"
$ws.Range("G25").Clear()

# Update the active selection to A25
$ws.Range("A25").Select()
